$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes all existing data rows down by one)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the game-reset / initial event entry
$ws.Cells.Item(2, 1).Value = "Init"
$ws.Cells.Item(2, 2).Value = "EV000"
$ws.Cells.Item(2, 3).Value = "Fade Out"

# The inserted row copied the bold header formatting; turn bold back off
# (Excel keeps the row's distinct font, just without the bold attribute)
$ws.Range("A2:C2").Font.Bold = $false

# Update the active selection to A2, matching the saved view state
$ws.Range("A2").Select()
